{"js": "// Update the date line and the 25 division-fact answer cells in the table.\n// Each cell is updated in place (range replace) so that the existing run\n// formatting (fonts / size) and paragraph formatting (alignment) carried by\n// the surrounding <w:r>/<w:pPr> is preserved instead of being stripped.\n\nfunction replaceRangeText(range, newText) {\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\n// 1. Title / date paragraph.\nconst dateResults = context.document.body.search(\"2024-05-04 Saturday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  replaceRangeText(dateResults.items[0], \"2024-05-05 Sunday\");\n  await context.sync();\n}\n\n// 2. Table of division problems.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, colIndex, newValue) \u2014 rowIndex is 0-based across ALL table\n// rows (including the blank spacer rows), matching the document's layout\n// where data only lives in rows 0, 4, 8, 12, 16.\nconst updates = [\n  [0, 0, \"74\u00f75=14, 4\"],\n  [0, 1, \"24\u00f72=12, 0\"],\n  [0, 2, \"61\u00f78=7, 5\"],\n  [0, 3, \"66\u00f75=13, 1\"],\n  [0, 4, \"85\u00f72=42, 1\"],\n\n  [4, 0, \"53\u00f75=10, 3\"],\n  [4, 1, \"77\u00f75=15, 2\"],\n  [4, 2, \"33\u00f76=5, 3\"],\n  [4, 3, \"12\u00f75=2, 2\"],\n  [4, 4, \"71\u00f79=7, 8\"],\n\n  [8, 0, \"62\u00f76=10, 2\"],\n  [8, 1, \"19\u00f72=9, 1\"],\n  [8, 2, \"79\u00f73=26, 1\"],\n  [8, 3, \"23\u00f79=2, 5\"],\n  [8, 4, \"73\u00f76=12, 1\"],\n\n  [12, 0, \"84\u00f73=28, 0\"],\n  [12, 1, \"66\u00f75=13, 1\"],\n  [12, 2, \"50\u00f73=16, 2\"],\n  [12, 3, \"66\u00f78=8, 2\"],\n  [12, 4, \"36\u00f79=4, 0\"],\n\n  [16, 0, \"88\u00f76=14, 4\"],\n  [16, 1, \"39\u00f73=13, 0\"],\n  // [16, 2] \"47\u00f77=6, 5\" is unchanged in the diff \u2014 skip it.\n  [16, 3, \"71\u00f72=35, 1\"],\n  [16, 4, \"32\u00f76=5, 2\"],\n];\n\nfor (const [rowIndex, colIndex, newValue] of updates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const range = cell.body.getRange();\n  replaceRangeText(range, newValue);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-fact answer cells in the table.\n# Each cell/paragraph Range.Text is assigned in place, which keeps the\n# existing run formatting (fonts, size) and paragraph formatting\n# (alignment) carried by the surrounding <w:rPr>/<w:pPr> untouched instead\n# of resetting it (unlike Range.Find/Execute, which in this host searches\n# past the supplied range's bounds instead of staying scoped to it).\n\n$d = $word.ActiveDocument\n\n# 1. Title / date paragraph.\n$d.Paragraphs(1).Range.Text = \"2024-05-05 Sunday\"\n\n# 2. Table of division problems.\n$t = $d.Tables.Item(1)\n\n# (rowIndex, colIndex, newValue) \u2014 Word COM Cell() is 1-based; data only\n# lives in rows 1, 5, 9, 13, 17 (the others are blank spacer rows).\n$updates = @(\n    @(1, 1, \"74\u00f75=14, 4\"),\n    @(1, 2, \"24\u00f72=12, 0\"),\n    @(1, 3, \"61\u00f78=7, 5\"),\n    @(1, 4, \"66\u00f75=13, 1\"),\n    @(1, 5, \"85\u00f72=42, 1\"),\n\n    @(5, 1, \"53\u00f75=10, 3\"),\n    @(5, 2, \"77\u00f75=15, 2\"),\n    @(5, 3, \"33\u00f76=5, 3\"),\n    @(5, 4, \"12\u00f75=2, 2\"),\n    @(5, 5, \"71\u00f79=7, 8\"),\n\n    @(9, 1, \"62\u00f76=10, 2\"),\n    @(9, 2, \"19\u00f72=9, 1\"),\n    @(9, 3, \"79\u00f73=26, 1\"),\n    @(9, 4, \"23\u00f79=2, 5\"),\n    @(9, 5, \"73\u00f76=12, 1\"),\n\n    @(13, 1, \"84\u00f73=28, 0\"),\n    @(13, 2, \"66\u00f75=13, 1\"),\n    @(13, 3, \"50\u00f73=16, 2\"),\n    @(13, 4, \"66\u00f78=8, 2\"),\n    @(13, 5, \"36\u00f79=4, 0\"),\n\n    @(17, 1, \"88\u00f76=14, 4\"),\n    @(17, 2, \"39\u00f73=13, 0\"),\n    # (17, 3) \"47\u00f77=6, 5\" is unchanged in the diff \u2014 skip it.\n    @(17, 4, \"71\u00f72=35, 1\"),\n    @(17, 5, \"32\u00f76=5, 2\")\n)\n\nforeach ($u in $updates) {\n    $rowIndex = $u[0]\n    $colIndex = $u[1]\n    $newText = $u[2]\n    $t.Cell($rowIndex, $colIndex).Range.Text = $newText\n}\n"}
